$d = $word.ActiveDocument
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Do zrobienia</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve">Z dnia </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="00B050"/>
        </w:rPr>
        <w:t>03-03-2015</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve"> 14:00 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t xml:space="preserve">wykonać </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve">do dnia </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>10-03-2015</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>1. Protokół założycielski</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> (PAWEŁ BIEŻUŃSKI) </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
        </w:rPr>
        <w:t>status</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>IN-PROGRESS</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>2. Specyfikacja wymagań</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>PATRYK TELUS</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve">status: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>IN-PROGRESS</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">3. </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Makiety</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">(PAWEŁ GOŁUŃSKI) </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">status: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="FF0000"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>IN-PROGRESS</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">4. </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Założenie</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> &quot;Issue Tracker&quot;</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> (ŁUKASZ FLAK) </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">status: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="FF0000"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>IN-PROGRESS</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>'
$d.Content.InsertXML($xml)
